# Applies the "feature selection" re-ranking edit described by the diff.
# Both worksheets (final_fail, final_gifted) keep the same row/column grid;
# only the feature label in column A (and, for a handful of rows whose
# ranking tier had ties, the RFE/RFECV/.../Total flags in columns B:J) are
# updated to reflect the re-sorted feature-importance table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "final_fail"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("final_fail")

$ws1.Range("A21").Value = "Number of days"
$ws1.Range("C21").Value = $true
$ws1.Range("D21").Value = $false

$ws1.Range("A22").Value = "Assignments viewed"
$ws1.Range("C22").Value = $false
$ws1.Range("D22").Value = $true

$ws1.Range("A25").Value = "Clicks on forum"

$ws1.Range("A26").Value = "Start of Session 5 (%)"
$ws1.Range("C26").Value = $true
$ws1.Range("D26").Value = $false

$ws1.Range("A27").Value = "Links viewed"
$ws1.Range("D27").Value = $true
$ws1.Range("J27").Value = 2

$ws1.Range("A28").Value = "Discussions viewed"
$ws1.Range("A29").Value = "Quizzes started"
$ws1.Range("A30").Value = "Forum posts"
$ws1.Range("A31").Value = "Assignments submitted"
$ws1.Range("A32").Value = "Number of sessions"

# ---------------------------------------------------------------------
# Sheet "final_gifted"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("final_gifted")

$ws2.Range("A5").Value = "Average session duration (min)"
$ws2.Range("A6").Value = "Clicks (% of course total)"
$ws2.Range("A7").Value = "Clicks per session"
$ws2.Range("A8").Value = "Clicks per day"
$ws2.Range("A9").Value = "On/off campus click ratio"
$ws2.Range("A10").Value = "Clicks on course"
$ws2.Range("A11").Value = "Start of Session 1 (%)"
$ws2.Range("A12").Value = "Start of Session 2 (%)"
$ws2.Range("A13").Value = "Start of Session 3 (%)"
$ws2.Range("A14").Value = "Days with no interaction"

$ws2.Range("A16").Value = "Days with no interaction (%)"
$ws2.Range("F16").Value = $true
$ws2.Range("J16").Value = 5

$ws2.Range("A17").Value = "Average grade of assignments"
$ws2.Range("C17").Value = $true
$ws2.Range("F17").Value = $false

$ws2.Range("A19").Value = "Start of Session 4 (%)"

$ws2.Range("A20").Value = "Clicks on forum"
$ws2.Range("C20").Value = $false
$ws2.Range("D20").Value = $true

$ws2.Range("A21").Value = "Files downloaded"
$ws2.Range("B21").Value = $false
$ws2.Range("J21").Value = 2

$ws2.Range("A22").Value = "Forum posts"
$ws2.Range("A23").Value = "Clicks on folder"
$ws2.Range("A24").Value = "Number of sessions"

$ws2.Range("A25").Value = "Start of Session 6 (%)"
$ws2.Range("B25").Value = $true
$ws2.Range("D25").Value = $false

$ws2.Range("A26").Value = "Start of Session 5 (%)"
$ws2.Range("A27").Value = "Assignments viewed"
$ws2.Range("A28").Value = "Discussions viewed"
$ws2.Range("A29").Value = "Quizzes started"
$ws2.Range("A30").Value = "Assignments submitted"
$ws2.Range("A31").Value = "Start of Session 10 (%)"
$ws2.Range("A32").Value = "Number of days"
$ws2.Range("A33").Value = "Start of Session 9 (%)"
$ws2.Range("A34").Value = "Start of Session 8 (%)"
$ws2.Range("A35").Value = "Start of Session 7 (%)"
$ws2.Range("A36").Value = "Submissions (% of course total)"
$ws2.Range("A37").Value = "Clicks on forum"
